$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (20, 22, 22 stored width units)
$ws.Columns.Item(1).ColumnWidth = 19.17
$ws.Columns.Item(2).ColumnWidth = 21.17
$ws.Columns.Item(3).ColumnWidth = 21.17

# Update header row
$ws.Range("A1").Value = "var_1_input_object"
$ws.Range("B1").Value = "var_2_input_object_1"
$ws.Range("C1").Value = "var_3_input_object_2"

# Update data row
$ws.Range("A2").Value = "cynthiacole"
$ws.Range("B2").Value = "Sample text"
$ws.Range("C2").Value = "fphillips"
